$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / timestamp text update ---
$ws.Range("A1").Value = "Datos actualizados a 19 de Marzo de 2020 a las 23:45"

# --- Two new countries (Lituania, Ucrania) pulled up into the ranking,
#     which shifts the country names of the rows below them down by one
#     within their respective blocks (numbers stay attached to their row). ---
$ws.Range("A89").Value  = "Lituania"
$ws.Range("A90").Value  = "Oman"
$ws.Range("A91").Value  = "Estado de Palestina"
$ws.Range("A92").Value  = "Kazajistan"
$ws.Range("A93").Value  = "Azerbaiyan"
$ws.Range("A94").Value  = "Venezuela"

$ws.Range("A105").Value = "Ucrania"
$ws.Range("A106").Value = "Uzbekistan"
$ws.Range("A107").Value = "Martinica"
$ws.Range("A108").Value = "Afganistan"

# --- Updated statistics (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) ---

function Set-Row($row, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Range("B$row").Value = $b
    $ws.Range("C$row").Value = $c
    $ws.Range("D$row").Value = $d
    $ws.Range("E$row").Value = $e
    $ws.Range("F$row").Value = $f
    $ws.Range("G$row").Value = $g
    $ws.Range("H$row").Value = $h
}

# Estados Unidos
Set-Row 9 13737 4478 108 13428 64 51 201

# Canada
Set-Row 22 841 114 11 818 1 3 12

# Peru
Set-Row 49 234 89 1 232 7 1 1

# Rows 89-94 (Lituania, Oman, Estado de Palestina, Kazajistan, Azerbaiyan, Venezuela)
Set-Row 89 48 14 1 47 1 0 0
Set-Row 90 48 9 13 35 0 0 0
Set-Row 91 47 3 0 47 0 0 0
Set-Row 92 44 8 0 41 0 3 3
Set-Row 93 44 10 7 36 0 0 1
Set-Row 94 42 6 0 42 0 0 0

# Rows 105-108 (Ucrania, Uzbekistan, Martinica, Afganistan)
Set-Row 105 26 10 0 23 0 1 3
Set-Row 106 23 5 0 23 0 0 0
Set-Row 107 23 0 0 22 0 0 1
Set-Row 108 22 0 1 21 0 0 0
